$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

$ws_ALC.Cells.Item(97, 8).Value2 = 995.4286
$ws_ALC.Cells.Item(97, 10).Value2 = 1013.3333
$ws_ALC.Cells.Item(97, 12).Value2 = 3039.9999
$ws_ALC.Cells.Item(97, 14).Value2 = -4031.9999
$ws_ALC.Cells.Item(100, 8).Value2 = 2285.2856
$ws_ALC.Cells.Item(100, 9).Value2 = 1000
$ws_ALC.Cells.Item(100, 10).Value2 = 2799.4
$ws_ALC.Cells.Item(100, 11).Value2 = 1000
$ws_ALC.Cells.Item(100, 12).Value2 = 2799.4
$ws_ALC.Cells.Item(100, 13).Value2 = -459
$ws_ALC.Cells.Item(100, 14).Value2 = -3881.4
$ws_ALC.Cells.Item(112, 8).Value2 = 1711.4872
$ws_ALC.Cells.Item(112, 9).Value2 = 1149.75
$ws_ALC.Cells.Item(112, 10).Value2 = 1775.6857
$ws_ALC.Cells.Item(112, 11).Value2 = 3449.25
$ws_ALC.Cells.Item(112, 12).Value2 = 5327.0571
$ws_ALC.Cells.Item(112, 13).Value2 = -2341.25
$ws_ALC.Cells.Item(112, 14).Value2 = -7543.0571
$ws_ALC.Cells.Item(116, 8).Value2 = 9109.071
$ws_ALC.Cells.Item(116, 9).Value2 = 14377.75
$ws_ALC.Cells.Item(116, 10).Value2 = 2084.1667
$ws_ALC.Cells.Item(116, 11).Value2 = 14377.75
$ws_ALC.Cells.Item(116, 12).Value2 = 2084.1667
$ws_ALC.Cells.Item(116, 13).Value2 = -10935.75
$ws_ALC.Cells.Item(116, 14).Value2 = -8968.1667
$ws_ALC.Cells.Item(133, 8).Value2 = 59444.445
$ws_ALC.Cells.Item(133, 10).Value2 = 59444.445
$ws_ALC.Cells.Item(133, 12).Value2 = 59444.445
$ws_ALC.Cells.Item(133, 14).Value2 = -69564.44500000001
$ws_ALC.Cells.Item(136, 8).Value2 = 72058.42999999999
$ws_ALC.Cells.Item(136, 10).Value2 = 72058.42999999999
$ws_ALC.Cells.Item(136, 12).Value2 = 72058.42999999999
$ws_ALC.Cells.Item(136, 14).Value2 = -82258.42999999999
$ws_ARM.Cells.Item(32, 8).Value2 = 4140.703
$ws_ARM.Cells.Item(32, 9).Value2 = 3664.1147
$ws_ARM.Cells.Item(32, 10).Value2 = 13831.333
$ws_ARM.Cells.Item(32, 11).Value2 = 3664.1147
$ws_ARM.Cells.Item(32, 12).Value2 = 13831.333
$ws_ARM.Cells.Item(32, 13).Value2 = -3377.1147
$ws_ARM.Cells.Item(32, 14).Value2 = -14405.333
$ws_ARM.Cells.Item(45, 8).Value2 = 1641.8182
$ws_ARM.Cells.Item(45, 9).Value2 = 1123
$ws_ARM.Cells.Item(45, 10).Value2 = 1836.375
$ws_ARM.Cells.Item(45, 11).Value2 = 1123
$ws_ARM.Cells.Item(45, 12).Value2 = 1836.375
$ws_ARM.Cells.Item(45, 13).Value2 = -746
$ws_ARM.Cells.Item(45, 14).Value2 = -2590.375
$ws_ARM.Cells.Item(63, 8).Value2 = 1619.6
$ws_ARM.Cells.Item(63, 9).Value2 = 1400
$ws_ARM.Cells.Item(63, 11).Value2 = 1400
$ws_ARM.Cells.Item(63, 13).Value2 = -714
$ws_ARM.Cells.Item(66, 8).Value2 = 1619.6
$ws_ARM.Cells.Item(66, 9).Value2 = 1400
$ws_ARM.Cells.Item(66, 11).Value2 = 7000
$ws_ARM.Cells.Item(66, 13).Value2 = -3568
$ws_ARM.Cells.Item(110, 8).Value2 = 453.16666
$ws_ARM.Cells.Item(110, 9).Value2 = 453.16666
$ws_ARM.Cells.Item(110, 11).Value2 = 453.16666
$ws_ARM.Cells.Item(110, 13).Value2 = 1591.83334
$ws_ARM.Cells.Item(132, 8).Value2 = 1871.3182
$ws_ARM.Cells.Item(132, 9).Value2 = 1375.5294
$ws_ARM.Cells.Item(132, 10).Value2 = 3557
$ws_ARM.Cells.Item(132, 11).Value2 = 4126.5882
$ws_ARM.Cells.Item(132, 12).Value2 = 10671
$ws_ARM.Cells.Item(132, 13).Value2 = -1596.5882
$ws_ARM.Cells.Item(132, 14).Value2 = -15731
$ws_BSM.Cells.Item(86, 8).Value2 = 108580.31
$ws_BSM.Cells.Item(86, 9).Value2 = 3532.875
$ws_BSM.Cells.Item(86, 10).Value2 = 668833.3
$ws_BSM.Cells.Item(86, 11).Value2 = 3532.875
$ws_BSM.Cells.Item(86, 12).Value2 = 668833.3
$ws_BSM.Cells.Item(86, 13).Value2 = -2409.875
$ws_BSM.Cells.Item(86, 14).Value2 = -671079.3
$ws_BSM.Cells.Item(89, 8).Value2 = 108580.31
$ws_BSM.Cells.Item(89, 9).Value2 = 3532.875
$ws_BSM.Cells.Item(89, 10).Value2 = 668833.3
$ws_BSM.Cells.Item(89, 11).Value2 = 17664.375
$ws_BSM.Cells.Item(89, 12).Value2 = 3344166.5
$ws_BSM.Cells.Item(89, 13).Value2 = -12048.375
$ws_BSM.Cells.Item(89, 14).Value2 = -3355398.5
$ws_BSM.Cells.Item(94, 8).Value2 = 2446.7
$ws_BSM.Cells.Item(94, 9).Value2 = 746
$ws_BSM.Cells.Item(94, 11).Value2 = 746
$ws_BSM.Cells.Item(94, 13).Value2 = -295
$ws_BSM.Cells.Item(99, 8).Value2 = 1562.375
$ws_BSM.Cells.Item(99, 9).Value2 = 1300.8
$ws_BSM.Cells.Item(99, 10).Value2 = 1998.3334
$ws_BSM.Cells.Item(99, 11).Value2 = 1300.8
$ws_BSM.Cells.Item(99, 12).Value2 = 1998.3334
$ws_BSM.Cells.Item(99, 13).Value2 = 197.2
$ws_BSM.Cells.Item(99, 14).Value2 = -4994.3334
$ws_BSM.Cells.Item(105, 8).Value2 = 2250.3438
$ws_BSM.Cells.Item(105, 9).Value2 = 2103.862
$ws_BSM.Cells.Item(105, 11).Value2 = 2103.862
$ws_BSM.Cells.Item(105, 13).Value2 = -356.8620000000001
$ws_BSM.Cells.Item(127, 8).Value2 = 31073.334
$ws_BSM.Cells.Item(127, 10).Value2 = 31073.334
$ws_BSM.Cells.Item(127, 12).Value2 = 31073.334
$ws_BSM.Cells.Item(127, 14).Value2 = -40993.334
$ws_CRP.Cells.Item(16, 8).Value2 = 586.15
$ws_CRP.Cells.Item(16, 9).Value2 = 519.9286
$ws_CRP.Cells.Item(16, 10).Value2 = 740.6667
$ws_CRP.Cells.Item(16, 11).Value2 = 519.9286
$ws_CRP.Cells.Item(16, 12).Value2 = 740.6667
$ws_CRP.Cells.Item(16, 13).Value2 = -232.9286
$ws_CRP.Cells.Item(16, 14).Value2 = -1314.6667
$ws_CRP.Cells.Item(31, 8).Value2 = 2798.3635
$ws_CRP.Cells.Item(31, 9).Value2 = 2500
$ws_CRP.Cells.Item(31, 10).Value2 = 2968.8572
$ws_CRP.Cells.Item(31, 11).Value2 = 2500
$ws_CRP.Cells.Item(31, 12).Value2 = 2968.8572
$ws_CRP.Cells.Item(31, 13).Value2 = -2205
$ws_CRP.Cells.Item(31, 14).Value2 = -3558.8572
$ws_CRP.Cells.Item(34, 8).Value2 = 2798.3635
$ws_CRP.Cells.Item(34, 9).Value2 = 2500
$ws_CRP.Cells.Item(34, 10).Value2 = 2968.8572
$ws_CRP.Cells.Item(34, 11).Value2 = 2500
$ws_CRP.Cells.Item(34, 12).Value2 = 2968.8572
$ws_CRP.Cells.Item(34, 13).Value2 = -2298
$ws_CRP.Cells.Item(34, 14).Value2 = -3372.8572
$ws_CRP.Cells.Item(58, 8).Value2 = 2072451.4
$ws_CRP.Cells.Item(58, 9).Value2 = 2175874
$ws_CRP.Cells.Item(58, 10).Value2 = 4000
$ws_CRP.Cells.Item(58, 11).Value2 = 2175874
$ws_CRP.Cells.Item(58, 12).Value2 = 4000
$ws_CRP.Cells.Item(58, 13).Value2 = -2175671
$ws_CRP.Cells.Item(58, 14).Value2 = -4406
$ws_CRP.Cells.Item(98, 8).Value2 = 50000
$ws_CRP.Cells.Item(98, 10).Value2 = 50000
$ws_CRP.Cells.Item(98, 12).Value2 = 50000
$ws_CRP.Cells.Item(98, 14).Value2 = -54492
$ws_CRP.Cells.Item(105, 8).Value2 = 877.25
$ws_CRP.Cells.Item(105, 9).Value2 = 861.7
$ws_CRP.Cells.Item(105, 10).Value2 = 955
$ws_CRP.Cells.Item(105, 11).Value2 = 861.7
$ws_CRP.Cells.Item(105, 12).Value2 = 955
$ws_CRP.Cells.Item(105, 13).Value2 = 885.3
$ws_CRP.Cells.Item(105, 14).Value2 = -4449
$ws_CRP.Cells.Item(113, 8).Value2 = 586.15
$ws_CRP.Cells.Item(113, 9).Value2 = 519.9286
$ws_CRP.Cells.Item(113, 10).Value2 = 740.6667
$ws_CRP.Cells.Item(113, 11).Value2 = 519.9286
$ws_CRP.Cells.Item(113, 12).Value2 = 740.6667
$ws_CRP.Cells.Item(113, 13).Value2 = 1650.0714
$ws_CRP.Cells.Item(113, 14).Value2 = -5080.6667
$ws_CRP.Cells.Item(122, 8).Value2 = 3445.7778
$ws_CRP.Cells.Item(122, 9).Value2 = 1833.1666
$ws_CRP.Cells.Item(122, 11).Value2 = 5499.4998
$ws_CRP.Cells.Item(122, 13).Value2 = -3049.4998
$ws_CRP.Cells.Item(136, 8).Value2 = 2072451.4
$ws_CRP.Cells.Item(136, 9).Value2 = 2175874
$ws_CRP.Cells.Item(136, 10).Value2 = 4000
$ws_CRP.Cells.Item(136, 11).Value2 = 6527622
$ws_CRP.Cells.Item(136, 12).Value2 = 12000
$ws_CRP.Cells.Item(136, 13).Value2 = -6525072
$ws_CRP.Cells.Item(136, 14).Value2 = -17100
$ws_CUL.Cells.Item(116, 8).Value2 = 55558100
$ws_CUL.Cells.Item(116, 9).Value2 = 97
$ws_CUL.Cells.Item(116, 10).Value2 = 62502850
$ws_CUL.Cells.Item(116, 11).Value2 = 291
$ws_CUL.Cells.Item(116, 12).Value2 = 187508550
$ws_CUL.Cells.Item(116, 13).Value2 = 3151
$ws_CUL.Cells.Item(116, 14).Value2 = -187515434
$ws_CUL.Cells.Item(122, 8).Value2 = 763.2
$ws_CUL.Cells.Item(122, 10).Value2 = 883.8
$ws_CUL.Cells.Item(122, 12).Value2 = 7954.2
$ws_CUL.Cells.Item(122, 14).Value2 = -12854.2
$ws_CUL.Cells.Item(134, 8).Value2 = 3189.28
$ws_CUL.Cells.Item(134, 9).Value2 = 2125.2
$ws_CUL.Cells.Item(134, 10).Value2 = 3898.6667
$ws_CUL.Cells.Item(134, 11).Value2 = 6375.599999999999
$ws_CUL.Cells.Item(134, 12).Value2 = 11696.0001
$ws_CUL.Cells.Item(134, 13).Value2 = -1305.599999999999
$ws_CUL.Cells.Item(134, 14).Value2 = -21836.0001
$ws_GSM.Cells.Item(39, 8).Value2 = 0
$ws_GSM.Cells.Item(39, 10).Value2 = 0
$ws_GSM.Cells.Item(39, 12).Value2 = 0
$ws_GSM.Cells.Item(39, 14).ClearContents()
$ws_GSM.Cells.Item(97, 8).Value2 = 797.9429
$ws_GSM.Cells.Item(97, 9).Value2 = 909.13635
$ws_GSM.Cells.Item(97, 10).Value2 = 609.7692
$ws_GSM.Cells.Item(97, 11).Value2 = 909.13635
$ws_GSM.Cells.Item(97, 12).Value2 = 609.7692
$ws_GSM.Cells.Item(97, 13).Value2 = -413.13635
$ws_GSM.Cells.Item(97, 14).Value2 = -1601.7692
$ws_GSM.Cells.Item(113, 8).Value2 = 876.2143
$ws_GSM.Cells.Item(113, 9).Value2 = 531.75
$ws_GSM.Cells.Item(113, 10).Value2 = 1335.5
$ws_GSM.Cells.Item(113, 11).Value2 = 531.75
$ws_GSM.Cells.Item(113, 12).Value2 = 1335.5
$ws_GSM.Cells.Item(113, 13).Value2 = 1638.25
$ws_GSM.Cells.Item(113, 14).Value2 = -5675.5
$ws_LTW.Cells.Item(68, 8).Value2 = 1563.6666
$ws_LTW.Cells.Item(68, 9).Value2 = 1563.6666
$ws_LTW.Cells.Item(68, 10).Value2 = 0
$ws_LTW.Cells.Item(68, 11).Value2 = 1563.6666
$ws_LTW.Cells.Item(68, 12).Value2 = 0
$ws_LTW.Cells.Item(68, 13).Value2 = -814.6666
$ws_LTW.Cells.Item(68, 14).ClearContents()
$ws_LTW.Cells.Item(71, 8).Value2 = 1563.6666
$ws_LTW.Cells.Item(71, 9).Value2 = 1563.6666
$ws_LTW.Cells.Item(71, 10).Value2 = 0
$ws_LTW.Cells.Item(71, 11).Value2 = 7818.333000000001
$ws_LTW.Cells.Item(71, 12).Value2 = 0
$ws_LTW.Cells.Item(71, 13).Value2 = -4074.333000000001
$ws_LTW.Cells.Item(71, 14).ClearContents()
$ws_LTW.Cells.Item(93, 8).Value2 = 16667224
$ws_LTW.Cells.Item(93, 9).Value2 = 573.2143
$ws_LTW.Cells.Item(93, 10).Value2 = 55556076
$ws_LTW.Cells.Item(93, 11).Value2 = 573.2143
$ws_LTW.Cells.Item(93, 12).Value2 = 55556076
$ws_LTW.Cells.Item(93, 13).Value2 = 674.7857
$ws_LTW.Cells.Item(93, 14).Value2 = -55558572
$ws_LTW.Cells.Item(124, 8).Value2 = 0
$ws_LTW.Cells.Item(124, 10).Value2 = 0
$ws_LTW.Cells.Item(124, 12).Value2 = 0
$ws_LTW.Cells.Item(124, 14).ClearContents()
$ws_WVR.Cells.Item(81, 8).Value2 = 1299
$ws_WVR.Cells.Item(81, 9).Value2 = 1332
$ws_WVR.Cells.Item(81, 10).Value2 = 1200
$ws_WVR.Cells.Item(81, 11).Value2 = 2664
$ws_WVR.Cells.Item(81, 12).Value2 = 2400
$ws_WVR.Cells.Item(81, 13).Value2 = -1603
$ws_WVR.Cells.Item(81, 14).Value2 = -4522
$ws_WVR.Cells.Item(84, 8).Value2 = 1299
$ws_WVR.Cells.Item(84, 9).Value2 = 1332
$ws_WVR.Cells.Item(84, 10).Value2 = 1200
$ws_WVR.Cells.Item(84, 11).Value2 = 13320
$ws_WVR.Cells.Item(84, 12).Value2 = 12000
$ws_WVR.Cells.Item(84, 13).Value2 = -8016
$ws_WVR.Cells.Item(84, 14).Value2 = -22608
$ws_WVR.Cells.Item(113, 8).Value2 = 503.73685
$ws_WVR.Cells.Item(113, 9).Value2 = 348.0909
$ws_WVR.Cells.Item(113, 10).Value2 = 717.75
$ws_WVR.Cells.Item(113, 11).Value2 = 1044.2727
$ws_WVR.Cells.Item(113, 12).Value2 = 2153.25
$ws_WVR.Cells.Item(113, 13).Value2 = 1125.7273
$ws_WVR.Cells.Item(113, 14).Value2 = -6493.25
$ws_WVR.Cells.Item(126, 8).Value2 = 1499.7693
$ws_WVR.Cells.Item(126, 9).Value2 = 1235.9412
$ws_WVR.Cells.Item(126, 11).Value2 = 3707.8236
$ws_WVR.Cells.Item(126, 13).Value2 = -1237.8236
